# Regen sval data to filter save games
# Update row 2 values (B2:E2 and G2) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1190320826869504
$ws.Range("C2").Value = 0.04071648406533734
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.406728370586922
